$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Helper: write a text value into a cell while keeping it as TEXT (string) type
# even when the value looks like a number (e.g. "6", "10", "53"). A leading
# apostrophe forces text entry (same as typing '6 in Excel); re-applying the
# "Normal" style afterwards strips the quote-prefix formatting that Excel
# would otherwise remember, leaving a plain text cell.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# New order rows 52-61 appended to the Orders sheet
Set-TextValue $ws.Range("C52") "315_尤加利叶圆叶_Eucalyptus Populus_undefined_1bunch"
Set-TextValue $ws.Range("F52") "6"

Set-TextValue $ws.Range("C53") "474_掌_anthurium_undefined_1bunch"
Set-TextValue $ws.Range("F53") "10"

Set-TextValue $ws.Range("C54") "300_白星_White Gypso_ gypsophila_1kg"
Set-TextValue $ws.Range("F54") "10"

Set-TextValue $ws.Range("C55") "728_金边朱蕉_undefined_undefined_1bunch"
Set-TextValue $ws.Range("F55") "8"

Set-TextValue $ws.Range("C56") "436_木百合_leucadendron _undefined_1bunch"
Set-TextValue $ws.Range("F56") "10"

Set-TextValue $ws.Range("C57") "496_大飞燕深蓝色_delphinium dark blue_undefined_1bunch"
Set-TextValue $ws.Range("F57") "5"

Set-TextValue $ws.Range("C58") "401_大飞燕白色_delphinium white_undefined_1bunch"
Set-TextValue $ws.Range("F58") "53"

Set-TextValue $ws.Range("A59") "9"
Set-TextValue $ws.Range("C59") "614_康乃馨绿_green_undefined_20stems"
Set-TextValue $ws.Range("F59") "5"

Set-TextValue $ws.Range("C60") "607_康乃馨芙蓉_light orange_undefined_20stems"
Set-TextValue $ws.Range("F60") "5"

Set-TextValue $ws.Range("C61") "615_康乃馨野马_horse_undefined_20stems"

# Summary sheet: G2 gains the newly-appended packages' numbers
$ws2 = $wb.Worksheets.Item("Summary")
Set-TextValue $ws2.Range("G2") "01013673102815383151019251841181010410197812530252525154060506101312251525655151210561010810553550"
